$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table (A2:E55) was sorted by "Reporting date" (column D) ascending
# (previously it was sorted descending, newest filing first -> oldest first).
$keyRange = $ws.Range("D1:D55")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange, 0, 1) | Out-Null
$ws.Sort.SetRange($ws.Range("A1:E55"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# The hyperlink that used to sit on the last (oldest) row now belongs to the
# first data row, since that row now holds the oldest filing.
$hyperlinkAddress = "https://www.sec.gov/Archives/edgar/data/1287750/000110465910055721/a10-17362_110q.htm"
$ws.Range("E55").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), $hyperlinkAddress) | Out-Null

# Update the frozen-pane view: scroll so the frozen pane's top-left visible
# cell is A2 and the active selection is E4.
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("E4").Select() | Out-Null
